$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.281.50'
$ws.Range('E2').Value = '  -3.42%  '
$ws.Range('D3').Value = '1.810.57'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4218'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3559'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07157'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8466'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.14%  '
$ws.Range('D12').Value = '1.828.72'
$ws.Range('E12').Value = '  -4.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.336'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.371'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06901'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.35%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008808'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.89%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.68%  '
$ws.Range('D21').Value = '27.685.02'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.094'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = '2.117.02'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.965'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.080'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.722'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08915'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7429'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.946'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.475'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.108'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.15%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.070'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05205'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01905'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.763'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1642'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4994'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.305'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.221'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '105.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06413'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4587'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.601'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.72%  '
